# Fix bug: exceeded requests in google drive
# - bump the date in A1 by one day (45310 -> 45311)
# - correct the price list values in D29:D32

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45311

$ws.Range("D29").Value = 4423.125
$ws.Range("D30").Value = 6063.75
$ws.Range("D31").Value = 6825
$ws.Range("D32").Value = 8728.125
